$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder/rename header row: A1=subject_code, B1=subject_name, C1=sem, D1=branchId
$ws.Range("A1").Value = "subject_code"
$ws.Range("B1").Value = "subject_name"
$ws.Range("C1").Value = "sem"
$ws.Range("D1").Value = "branchId"

# Move selection to H4 as in the target sheet view
$ws.Range("H4").Select()
